# "Sankey All Elc Sector fix"
#
# - Rename sheet "TS_Defs_Sankey" -> "Sankey_def" (the _xlnm._FilterDatabase
#   defined name tracks the sheet by index and gets its sheet-qualified
#   reference text updated automatically by the rename).
# - C3 on that sheet: "Power" -> "ElectricitySector".
# - A18 on that sheet: "TS_Defs: snk_attr=Sankey_ResDetV0"
#       -> "~TS_Defs: snk_attr=Sankey_Test" (block header replaced/renamed).
# - Column C width narrowed from 96 to ~78.29 characters.
# - View: zoom to 120%, selection moved to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS_Defs_Sankey")

$ws.Name = "Sankey_def"

$ws.Range("C3").Value = "ElectricitySector"
$ws.Range("A18").Value = "~TS_Defs: snk_attr=Sankey_Test"

$ws.Columns.Item(3).ColumnWidth = 77.5

$ws.Activate()
[void]$ws.Range("C3").Select()
$excel.ActiveWindow.Zoom = 120
